$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new product row at row 4 by shifting existing rows (4..59) down
# to (5..60) via plain value copies - this avoids COM "Insert" row-shift
# semantics forking a brand-new cell style, keeping styles.xml untouched
# (cell styles stay bound to row position, matching the source diff).
for ($r = 59; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":P" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":P" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# Row 60 is brand new (previously the sheet ended at row 59), so its A/O
# cells came out with no explicit style. Copy just the cell format from the
# old last row (still intact at row 59) so A60/O60 keep the same style
# indices (1 / 3) as every other data row, without minting a new cellXfs.
$ws.Cells.Item(59, 1).Copy() | Out-Null
$ws.Cells.Item(60, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(59, 15).Copy() | Out-Null
$ws.Cells.Item(60, 15).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the freed-up row 4 with the new "Gaseosa Pepsi con azúcar" record.
$ws.Cells.Item(4, 1).Value2 = 7791813888406
$ws.Cells.Item(4, 2).Value2 = "Gaseosa"
$ws.Cells.Item(4, 3).Value2 = "original"
$ws.Cells.Item(4, 4).Value2 = "con azúcar"
$ws.Cells.Item(4, 5).Value2 = "Pepsi"
$ws.Cells.Item(4, 6).Value2 = 500
$ws.Cells.Item(4, 7).Value2 = "ml."
$ws.Cells.Item(4, 8).Value2 = "botella"
$ws.Cells.Item(4, 9).Value2 = "Gaseosas"
$ws.Cells.Item(4, 10).Value2 = "Argentina"
$ws.Cells.Item(4, 11).Value2 = 6
$ws.Cells.Item(4, 12).Value2 = $false
$ws.Cells.Item(4, 13).Value2 = $true
$ws.Cells.Item(4, 14).Value2 = "C:\VentaSoft\Imágenes de artículos\7791813888406.png"
$ws.Cells.Item(4, 15).Value2 = $true
$ws.Cells.Item(4, 16).Value2 = $true
